$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A, B, C across rows 2-21 (Top20 hot stock lists).
$data = @(
    @("三花智控", "三花智控", "常山北明"),
    @("常山北明", "闻泰科技", "白银有色"),
    @("山子高科", "白银有色", "楚江新材"),
    @("四方股份", "常山北明", "山子高科"),
    @("白银有色", "山子高科", "三花智控"),
    @("中兴通讯", "中兴通讯", "华建集团"),
    @("楚江新材", "四方股份", "合肥城建"),
    @("合肥城建", "楚江新材", "海马汽车"),
    @("闻泰科技", "通富微电", "大有能源"),
    @("海马汽车", "香农芯创", "拓维信息"),
    @("拓维信息", "拓维信息", "安泰科技"),
    @("香农芯创", "合肥城建", "天普股份"),
    @("大有能源", "大有能源", "合锻智能"),
    @("宝泰隆", "宝泰隆", "闻泰科技"),
    @("华建集团", "中国核建", "中兴通讯"),
    @("阳光电源", "海马汽车", "宝泰隆"),
    @("德明利", "阳光电源", "天际股份"),
    @("通富微电", "北方稀土", "阳光电源"),
    @("中信证券", "通达股份", "黄河旋风"),
    @("京泉华", "成飞集成", "华胜天成")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
